$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily record as row 75 (dimension grows from A1:D74 to A1:D75).
# The date column stores plain text like "2025/10/07" (same as every other
# row), so a leading apostrophe forces text entry instead of Excel's
# automatic date-value conversion; resetting the style back to "Normal"
# afterwards strips the quote-prefix formatting Excel applied, leaving a
# plain unstyled text cell just like the rest of the column.
$ws.Range("A75").Value = "'2025/10/07"
$ws.Range("A75").Style = "Normal"

$ws.Range("B75").Value = "火"
$ws.Range("C75").Value = 22
$ws.Range("D75").Value = 11
